$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 1.02
    "C2" = 1.063698992546387
    "D2" = 1.045782457955618
    "E2" = 1.068376116009812
    "F2" = 1.077582912093766
    "I2" = 1.044904775299537
    "J2" = 1.068663028303149
    "K2" = 1.048549731816482
    "L2" = 1.071081464442082
    "M2" = 1.080263831214167
    "N2" = 1.070180652380701
    "B3" = 1.02
    "C3" = 1.06529527807644
    "D3" = 1.046501588922199
    "E3" = 1.069841698826519
    "F3" = 1.079231884031212
    "I3" = 1.045274168214591
    "J3" = 1.069911169654141
    "K3" = 1.049080613309382
    "L3" = 1.072361399854886
    "M3" = 1.081728496877777
    "N3" = 1.071430566235576
    "B4" = 1.02
    "C4" = 1.066326355271698
    "D4" = 1.04696614021125
    "E4" = 1.070788577714938
    "F4" = 1.080297662927486
    "I4" = 1.045511283590716
    "J4" = 1.070716587263323
    "K4" = 1.049422685365117
    "L4" = 1.073187636220298
    "M4" = 1.082674534170647
    "N4" = 1.072237127630146
    "B5" = 1.02
    "C5" = 1.066759392608112
    "D5" = 1.047161253765016
    "E5" = 1.07118630629387
    "F5" = 1.080745434446511
    "I5" = 1.045610512600977
    "J5" = 1.071054662424139
    "K5" = 1.049566148702242
    "L5" = 1.073534522472364
    "M5" = 1.083071849868183
    "N5" = 1.072575682896468
    "B6" = 1.02
    "C6" = 1.066832076693286
    "D6" = 1.047194003457137
    "E6" = 1.071253067010076
    "F6" = 1.080820600989652
    "I6" = 1.045627147010137
    "J6" = 1.071111396313383
    "K6" = 1.049590216724245
    "L6" = 1.073592739236417
    "M6" = 1.083138537805171
    "N6" = 1.072632497354343
    "B7" = 1.02
    "C7" = 1.066332143205659
    "D7" = 1.046968748048384
    "E7" = 1.070793893506905
    "F7" = 1.080303647168824
    "I7" = 1.045512611276049
    "J7" = 1.070721106683654
    "K7" = 1.049424603677668
    "L7" = 1.073192273143444
    "M7" = 1.082679844677392
    "N7" = 1.072241653468571
    "B8" = 1.02
    "C8" = 1.064238848714735
    "D8" = 1.046025652985839
    "E8" = 1.068871721962674
    "F8" = 1.078140447279862
    "I8" = 1.045030009417458
    "J8" = 1.069085306305535
    "K8" = 1.048729445786476
    "L8" = 1.071514436501204
    "M8" = 1.080759178787896
    "N8" = 1.070603530066287
    "B9" = 1.02
    "C9" = 1.060535773092363
    "D9" = 1.044357793157307
    "E9" = 1.065473130716388
    "F9" = 1.074318882725427
    "I9" = 1.044164905025057
    "J9" = 1.066185529719609
    "K9" = 1.047493345634565
    "L9" = 1.068542468707098
    "M9" = 1.077361325499974
    "N9" = 1.067699635465
    "B10" = 1.02
    "C10" = 1.058056746553935
    "D10" = 1.043241742116987
    "E10" = 1.063199172374968
    "F10" = 1.071764048198056
    "I10" = 1.043578150088268
    "J10" = 1.064240250066305
    "K10" = 1.046661664777704
    "L10" = 1.066550325600353
    "M10" = 1.075086550884299
    "N10" = 1.065751593291461
    "B11" = 1.02
    "C11" = 1.056980718591086
    "D11" = 1.042757471987083
    "E11" = 1.062212456264903
    "F11" = 1.070655956343567
    "I11" = 1.043321671594567
    "J11" = 1.063394947884082
    "K11" = 1.04629970392311
    "L11" = 1.06568503054958
    "M11" = 1.074099166464826
    "N11" = 1.064905090683182
    "B12" = 1.02
    "C12" = 1.056580632442411
    "D12" = 1.042577438128401
    "E12" = 1.061845623929459
    "F12" = 1.070244075198306
    "I12" = 1.043226039304036
    "J12" = 1.063080507614276
    "K12" = 1.046164976984423
    "L12" = 1.065363208986827
    "M12" = 1.073732038193133
    "N12" = 1.064590203872127
    "B13" = 1.02
    "C13" = 1.056666470637655
    "D13" = 1.04261606302398
    "E13" = 1.061924325444618
    "F13" = 1.070332438249143
    "I13" = 1.043246569305081
    "J13" = 1.063147976966571
    "K13" = 1.046193889015889
    "L13" = 1.065432259617818
    "M13" = 1.073810805352313
    "N13" = 1.064657769038641
    "B14" = 1.02
    "C14" = 1.056947655603769
    "D14" = 1.042742593494295
    "E14" = 1.06218214041711
    "F14" = 1.070621916029748
    "I14" = 1.043313774056409
    "J14" = 1.063368965537682
    "K14" = 1.046288573049404
    "L14" = 1.065658437146381
    "M14" = 1.074068827130123
    "N14" = 1.064879071438871
    "B15" = 1.02
    "C15" = 1.05712084938845
    "D15" = 1.042820532570218
    "E15" = 1.062340945757569
    "F15" = 1.07080023454129
    "I15" = 1.043355132705791
    "J15" = 1.063505063002957
    "K15" = 1.04634687402037
    "L15" = 1.065797737655129
    "M15" = 1.074227753558139
    "N15" = 1.065015362178157
    "B16" = 1.02
    "C16" = 1.058128103341524
    "D16" = 1.043273860005507
    "E16" = 1.063264612861415
    "F16" = 1.071837549051351
    "I16" = 1.043595120719093
    "J16" = 1.064296286350891
    "K16" = 1.046685648000637
    "L16" = 1.066607695050289
    "M16" = 1.075152029057105
    "N16" = 1.065807709153999
    "B17" = 1.02
    "C17" = 1.058759224539016
    "D17" = 1.043557947647393
    "E17" = 1.063843441710442
    "F17" = 1.072487730701434
    "I17" = 1.043745011791451
    "J17" = 1.064791794680262
    "K17" = 1.046897658255089
    "L17" = 1.067115034873515
    "M17" = 1.075731154963705
    "N17" = 1.066303921162033
    "B18" = 1.02
    "C18" = 1.059127097465181
    "D18" = 1.043723553608367
    "E18" = 1.06418086304945
    "F18" = 1.072866794699855
    "I18" = 1.043832208399309
    "J18" = 1.065080529090729
    "K18" = 1.047021143122526
    "L18" = 1.067410699275175
    "M18" = 1.076068719218816
    "N18" = 1.066593065608482
    "B19" = 1.02
    "C19" = 1.059252490682643
    "D19" = 1.043780004481103
    "E19" = 1.064295881427802
    "F19" = 1.07299601622236
    "I19" = 1.043861900865765
    "J19" = 1.065178931623482
    "K19" = 1.047063218301955
    "L19" = 1.067511469655176
    "M19" = 1.076183781191399
    "N19" = 1.066691607884119
    "B20" = 1.02
    "C20" = 1.058691537088907
    "D20" = 1.043527477838675
    "E20" = 1.063781359584692
    "F20" = 1.072417990645634
    "I20" = 1.043728953947482
    "J20" = 1.064738661100528
    "K20" = 1.046874929900369
    "L20" = 1.067060628899316
    "M20" = 1.075669044108831
    "N20" = 1.066250712126523
    "B21" = 1.02
    "C21" = 1.056864864817488
    "D21" = 1.042705337703933
    "E21" = 1.062106229303943
    "F21" = 1.070536680038315
    "I21" = 1.043293994023297
    "J21" = 1.063303902612335
    "K21" = 1.046260698671267
    "L21" = 1.065591844961754
    "M21" = 1.073992856437763
    "N21" = 1.064813916116707
    "B22" = 1.02
    "C22" = 1.05571403482316
    "D22" = 1.04218753130085
    "E22" = 1.061051140191156
    "F22" = 1.069352162159003
    "I22" = 1.043018405707533
    "J22" = 1.062399161999282
    "K22" = 1.045872894194351
    "L22" = 1.064665971247322
    "M22" = 1.072936826348854
    "N22" = 1.063907890668208
    "B23" = 1.02
    "C23" = 1.05632433629355
    "D23" = 1.042462115791994
    "E23" = 1.061610643329847
    "F23" = 1.06998025895423
    "I23" = 1.043164701373419
    "J23" = 1.062879036494249
    "K23" = 1.04607863041073
    "L23" = 1.065157024156174
    "M23" = 1.07349685462464
    "N23" = 1.06438844664
    "B24" = 1.02
    "C24" = 1.058722122904733
    "D24" = 1.04354124613539
    "E24" = 1.063809412442341
    "F24" = 1.072449503718762
    "I24" = 1.043736210513908
    "J24" = 1.06476267076569
    "K24" = 1.046885200413026
    "L24" = 1.067085213415532
    "M24" = 1.075697110041701
    "N24" = 1.066274755888162
    "B25" = 1.02
    "C25" = 1.061494878847631
    "D25" = 1.044789697470608
    "E25" = 1.066353160913011
    "F25" = 1.075308061544601
    "I25" = 1.044390311022184
    "J25" = 1.066937288484681
    "K25" = 1.047814240137005
    "L25" = 1.069312667739351
    "M25" = 1.078241393988692
    "N25" = 1.068452461813747
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

Write-Host "Updated $($updates.Count) cells"